$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed consistency/opportunity/percentage figures for the 2015/2016 IN re-pull
# Each tuple is (row, consistency, opportunity, percentage)
$data = @(
    @(2, 7, 9, 0.7777777777777778),
    @(4, 11, 11, 1),
    @(5, 24, 25, 0.96),
    @(6, 8, 9, 0.8888888888888888),
    @(7, 14, 15, 0.9333333333333333),
    @(8, 10, 10, 1),
    @(9, 12, 12, 1),
    @(10, 8, 9, 0.8888888888888888),
    @(11, 6, 7, 0.8571428571428571),
    @(12, 23, 25, 0.92),
    @(13, 12, 13, 0.9230769230769231),
    @(14, 11, 11, 1),
    @(16, 21, 24, 0.875),
    @(17, 21, 22, 0.9545454545454546),
    @(18, 24, 26, 0.9230769230769231),
    @(19, 21, 23, 0.9130434782608695),
    @(20, 9, 9, 1),
    @(21, 14, 14, 1),
    @(22, 17, 22, 0.7727272727272727),
    @(23, 11, 13, 0.8461538461538461),
    @(24, 8, 11, 0.7272727272727273),
    @(25, 24, 24, 1),
    @(26, 7, 7, 1),
    @(27, 9, 9, 1),
    @(28, 15, 16, 0.9375),
    @(29, 19, 19, 1),
    @(30, 5, 6, 0.8333333333333334),
    @(31, 21, 23, 0.9130434782608695),
    @(32, 11, 11, 1),
    @(33, 5, 5, 1),
    @(34, 23, 24, 0.9583333333333334),
    @(35, 10, 10, 1),
    @(36, 8, 13, 0.6153846153846154),
    @(37, 12, 12, 1),
    @(38, 15, 16, 0.9375),
    @(40, 11, 11, 1),
    @(41, 24, 24, 1),
    @(43, 26, 30, 0.8666666666666667),
    @(44, 7, 7, 1),
    @(45, 11, 12, 0.9166666666666666),
    @(47, 10, 10, 1),
    @(48, 19, 25, 0.76),
    @(49, 12, 13, 0.9230769230769231),
    @(50, 18, 18, 1),
    @(51, 13, 14, 0.9285714285714286),
    @(53, 8, 10, 0.8),
    @(54, 16, 17, 0.9411764705882353),
    @(55, 25, 32, 0.78125),
    @(56, 12, 13, 0.9230769230769231),
    @(57, 5, 5, 1),
    @(58, 21, 25, 0.84),
    @(59, 16, 16, 1),
    @(60, 25, 25, 1),
    @(62, 8, 8, 1),
    @(63, 27, 27, 1),
    @(64, 21, 22, 0.9545454545454546),
    @(66, 7, 9, 0.7777777777777778),
    @(67, 9, 9, 1),
    @(69, 19, 20, 0.95),
    @(70, 17, 17, 1),
    @(71, 14, 20, 0.7),
    @(72, 3, 5, 0.6),
    @(73, 15, 19, 0.7894736842105263),
    @(74, 27, 32, 0.84375),
    @(75, 21, 21, 1),
    @(76, 4, 4, 1),
    @(77, 5, 5, 1),
    @(79, 26, 29, 0.896551724137931),
    @(80, 10, 11, 0.9090909090909091),
    @(81, 11, 11, 1),
    @(83, 11, 11, 1),
    @(85, 6, 6, 1),
    @(86, 23, 26, 0.8846153846153846),
    @(87, 15, 15, 1),
    @(89, 33, 33, 1),
    @(93, 12, 14, 0.8571428571428571),
    @(94, 22, 27, 0.8148148148148148),
    @(95, 13, 13, 1),
    @(96, 8, 10, 0.8),
    @(97, 19, 21, 0.9047619047619048),
    @(98, 10, 10, 1),
    @(99, 13, 14, 0.9285714285714286),
    @(100, 16, 16, 1),
    @(101, 21, 22, 0.9545454545454546)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    if ($null -ne $row[3]) {
        $ws.Range("D$r").Value = $row[3]
    }
}

Write-Host "Updated $($data.Count) rows in house_consistency_matrix"
